# Updates for Financials (GFI yearly financial statements)
# Refresh columns D:J (FY2017..FY2011) across Income Statement, Balance Sheet,
# and Cash Flow Statement sections with newly reported figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: Total Revenue
$ws.Range("D8").Value = 2577800
$ws.Range("E8").Value = 2761800
$ws.Range("F8").Value = 2666400
$ws.Range("G8").Value = 2545400
$ws.Range("H8").Value = 2868800
$ws.Range("I8").Value = 2906300
$ws.Range("J8").Value = 3530600

# Row 9: Cost of Revenue
$ws.Range("D9").Value = 2043000
$ws.Range("E9").Value = 2105100
$ws.Range("F9").Value = 2001200
$ws.Range("G9").Value = 2066100
$ws.Range("H9").Value = 1808100
$ws.Range("I9").Value = 1819900
$ws.Range("J9").Value = 3525600

# Row 10: Gross Profit
$ws.Range("D10").Value = 534800
$ws.Range("E10").Value = 656700
$ws.Range("F10").Value = 665200
$ws.Range("G10").Value = 479300
$ws.Range("H10").Value = 1060700
$ws.Range("I10").Value = 1086400
$ws.Range("J10").Value = 5000

# Row 12: Research Development
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "NA"
$ws.Range("H12").Value = 36200
$ws.Range("I12").Value = 77900
$ws.Range("J12").Value = 135300

# Row 14: Non Recurring
$ws.Range("D14").Value = 581500
$ws.Range("E14").Value = 239600
$ws.Range("F14").Value = 88200
$ws.Range("G14").Value = 230400
$ws.Range("H14").Value = 107200
$ws.Range("I14").Value = 365300
$ws.Range("J14").Value = 146600

# Row 15: Others
$ws.Range("D15").Value = "NA"
$ws.Range("E15").Value = "NA"
$ws.Range("F15").Value = "NA"
$ws.Range("G15").Value = "NA"
$ws.Range("H15").Value = 677300
$ws.Range("I15").Value = 568500
$ws.Range("J15").Value = 925600

# Row 17: Total Operating Expenses
$ws.Range("D17").Value = 2624500
$ws.Range("E17").Value = 2344700
$ws.Range("F17").Value = 2089400
$ws.Range("G17").Value = 2296500
$ws.Range("H17").Value = 2758900
$ws.Range("I17").Value = 3029700
$ws.Range("J17").Value = 2792000

# Row 18: Operating Income or Loss
$ws.Range("D18").Value = -46700
$ws.Range("E18").Value = 417100
$ws.Range("F18").Value = 577000
$ws.Range("G18").Value = 248900
$ws.Range("H18").Value = 109900
$ws.Range("I18").Value = -123400
$ws.Range("J18").Value = 738600

# Row 20: Total Other Income/Expenses Net
$ws.Range("D20").Value = -289000
$ws.Range("E20").Value = -196400
$ws.Range("F20").Value = -152200
$ws.Range("G20").Value = -173200
$ws.Range("H20").Value = 79100
$ws.Range("I20").Value = 33300
$ws.Range("J20").Value = 16400

# Row 21: Earnings Before Interest And Taxes
$ws.Range("D21").Value = 412400
$ws.Range("E21").Value = 892100
$ws.Range("F21").Value = 1034700
$ws.Range("G21").Value = 753000
$ws.Range("H21").Value = 757500
$ws.Range("I21").Value = 335700
$ws.Range("J21").Value = "NA"

# Row 22: Interest Expense
$ws.Range("D22").Value = 75000
$ws.Range("E22").Value = 68300
$ws.Range("F22").Value = 67400
$ws.Range("G22").Value = 71200
$ws.Range("H22").Value = 80800
$ws.Range("I22").Value = 72400
$ws.Range("J22").Value = 42300

# Row 23: Income Before Tax
$ws.Range("D23").Value = -410700
$ws.Range("E23").Value = 152400
$ws.Range("F23").Value = 357400
$ws.Range("G23").Value = 4500
$ws.Range("H23").Value = 108200
$ws.Range("I23").Value = -162500
$ws.Range("J23").Value = 712700

# Row 24: Income Tax Expense
$ws.Range("D24").Value = -65900
$ws.Range("E24").Value = 173200
$ws.Range("F24").Value = 189500
$ws.Range("G24").Value = 247100
$ws.Range("H24").Value = 121600
$ws.Range("I24").Value = 105700
$ws.Range("J24").Value = 359400

# Row 26: Income After Tax
$ws.Range("D26").Value = -344800
$ws.Range("E26").Value = -20800
$ws.Range("F26").Value = 167900
$ws.Range("G26").Value = -242600
$ws.Range("H26").Value = -13400
$ws.Range("I26").Value = -268200
$ws.Range("J26").Value = 353300

# Row 27: Net Income From Continuing Ops
$ws.Range("D27").Value = -348200
$ws.Range("E27").Value = -31800
$ws.Range("F27").Value = 157000
$ws.Range("G27").Value = -242100
$ws.Range("H27").Value = -27200
$ws.Range("I27").Value = -268400
$ws.Range("J27").Value = 292100

# Row 29: Discontinued Operations
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 13100
$ws.Range("F29").Value = 1200
$ws.Range("G29").Value = "NA"
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 20500
$ws.Range("J29").Value = 747100

# Row 32: Other Items
$ws.Range("D32").Value = 289000
$ws.Range("E32").Value = 196400
$ws.Range("F32").Value = 152200
$ws.Range("G32").Value = 173200
$ws.Range("H32").Value = -79100
$ws.Range("I32").Value = -33300
$ws.Range("J32").Value = -16400

# Row 33: Net Income
$ws.Range("D33").Value = -348200
$ws.Range("E33").Value = -18700
$ws.Range("F33").Value = 158200
$ws.Range("G33").Value = -242100
$ws.Range("H33").Value = -27200
$ws.Range("I33").Value = -247900
$ws.Range("J33").Value = 1039200

# Row 35: Net Income Applicable To Common Shares
$ws.Range("D35").Value = -348200
$ws.Range("E35").Value = -18700
$ws.Range("F35").Value = 158200
$ws.Range("G35").Value = -242100
$ws.Range("H35").Value = -27200
$ws.Range("I35").Value = -247900
$ws.Range("J35").Value = 1039200

# Row 41: Cash And Cash Equivalents
$ws.Range("E41").Value = 526700

# Row 43: Net Receivables
$ws.Range("F43").Value = 128200

# Row 44: Inventory
$ws.Range("E44").Value = 329400
$ws.Range("F44").Value = 298200

# Row 45: Other Current Assets
$ws.Range("E45").Value = 76800

# Row 47: Long Term Investments
$ws.Range("F47").Value = 140000

# Row 48: Property Plant and Equipment
$ws.Range("E48").Value = 4524600
$ws.Range("F48").Value = 4312400

# Row 52: Other Assets
$ws.Range("F52").Value = 221900

# Row 57: Accounts Payable
$ws.Range("F57").Value = 155300

# Row 59: Other Current Liabilities
$ws.Range("F59").Value = 350100

# Row 62: Other Liabilities
$ws.Range("E62").Value = 773900
$ws.Range("F62").Value = 784000

# Row 81: Net Income
$ws.Range("D81").Value = -348200
$ws.Range("E81").Value = -18700
$ws.Range("F81").Value = 158200
$ws.Range("G81").Value = -242100
$ws.Range("H81").Value = -27200
$ws.Range("I81").Value = -247900
$ws.Range("J81").Value = 1039200

# Row 83: Depreciation
$ws.Range("F83").Value = 609900
$ws.Range("G83").Value = 677300
$ws.Range("J83").Value = "NA"

# Row 89: Total Cash Flow From Operating Activities
$ws.Range("G89").Value = 743800

# Row 91: Capital Expenditures
$ws.Range("F91").Value = -634100
$ws.Range("G91").Value = -480500

# Row 94: Total Cash Flows From Investing Activities
$ws.Range("G94").Value = -423900
$ws.Range("J94").Value = "NA"

# Row 100: Total Cash Flows From Financing Activities
$ws.Range("G100").Value = -168200
$ws.Range("J100").Value = "NA"

# Row 101: Effect Of Exchange Rate Changes
$ws.Range("J101").Value = "NA"
